$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Alternate form factor work: update PCB component position/rotation values
$ws.Range("B52").Value = 76.77

$ws.Range("C53").Value = -45.8

$ws.Range("B54").Value = 97.45
$ws.Range("C54").Value = -52.5
$ws.Range("E54").Value = 180

$ws.Range("C55").Value = -48.65

$ws.Range("B56").Value = 149.19999999999999
$ws.Range("C56").Value = -54.5

$ws.Range("C57").Value = -109.02

$ws.Range("C60").Value = -55.25

$ws.Range("C61").Value = -55.25

$ws.Range("B62").Value = 86.05

# Update the view: zoom level and scroll/selection position
$excel.ActiveWindow.Zoom = 175
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C57").Select()
